$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.055589666666667
$ws.Range("H2").Value = 6.166769
$ws.Range("I2").Value = 0.1150400298148962
$ws.Range("J2").Value = 0.1150400298148962
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 4.543446666666667
$ws.Range("N2").Value = 13.63034
$ws.Range("O2").Value = 0.305297842646339
$ws.Range("P2").Value = 0.305297842646339
$ws.Range("Q2").Value = 9.339462019051112
$ws.Range("R2").Value = 84.05515817146001
$ws.Range("S2").Value = 0.03512147292045831
$ws.Range("T2").Value = 0.03512147292045831

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.055589666666667
$ws.Range("H3").Value = 6.166769
$ws.Range("I3").Value = 0.1150400298148962
$ws.Range("J3").Value = 0.1150400298148962
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 4.938922000000001
$ws.Range("N3").Value = 14.816766
$ws.Range("O3").Value = 0.3318718898278125
$ws.Range("P3").Value = 0.3318718898278125
$ws.Range("Q3").Value = 10.15239702767267
$ws.Range("R3").Value = 91.37157324905401
$ws.Range("S3").Value = 0.03817855210051748
$ws.Range("T3").Value = 0.03817855210051748

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.055589666666667
$ws.Range("H4").Value = 6.166769
$ws.Range("I4").Value = 0.1150400298148962
$ws.Range("J4").Value = 0.1150400298148962
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 3.044074333333333
$ws.Range("N4").Value = 9.132223
$ws.Range("O4").Value = 0.2045472072204566
$ws.Range("P4").Value = 0.2045472072204566
$ws.Range("Q4").Value = 6.257367744165221
$ws.Range("R4").Value = 56.316309697487
$ws.Range("S4").Value = 0.02353111681719506
$ws.Range("T4").Value = 0.02353111681719506

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 2.055589666666667
$ws.Range("H5").Value = 6.166769
$ws.Range("I5").Value = 0.1150400298148962
$ws.Range("J5").Value = 0.1150400298148962
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 2.355570666666666
$ws.Range("N5").Value = 7.066712
$ws.Range("O5").Value = 0.1582830603053919
$ws.Range("P5").Value = 0.1582830603053919
$ws.Range("Q5").Value = 4.842086721503111
$ws.Range("R5").Value = 43.578780493528
$ws.Range("S5").Value = 0.01820888797672529
$ws.Range("T5").Value = 0.01820888797672529

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 8.497489666666667
$ws.Range("H6").Value = 25.492469
$ws.Range("I6").Value = 0.4755576856884563
$ws.Range("J6").Value = 0.4755576856884562
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 4.543446666666667
$ws.Range("N6").Value = 13.63034
$ws.Range("O6").Value = 0.305297842646339
$ws.Range("P6").Value = 0.305297842646339
$ws.Range("Q6").Value = 38.60789110105112
$ws.Range("R6").Value = 347.47101990946
$ws.Range("S6").Value = 0.1451867354945715
$ws.Range("T6").Value = 0.1451867354945715

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 8.497489666666667
$ws.Range("H7").Value = 25.492469
$ws.Range("I7").Value = 0.4755576856884563
$ws.Range("J7").Value = 0.4755576856884562
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 4.938922000000001
$ws.Range("N7").Value = 14.816766
$ws.Range("O7").Value = 0.3318718898278125
$ws.Range("P7").Value = 0.3318718898278125
$ws.Range("Q7").Value = 41.96843865947267
$ws.Range("R7").Value = 377.715947935254
$ws.Range("S7").Value = 0.1578242278715689
$ws.Range("T7").Value = 0.1578242278715689

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 8.497489666666667
$ws.Range("H8").Value = 25.492469
$ws.Range("I8").Value = 0.4755576856884563
$ws.Range("J8").Value = 0.4755576856884562
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 3.044074333333333
$ws.Range("N8").Value = 9.132223
$ws.Range("O8").Value = 0.2045472072204566
$ws.Range("P8").Value = 0.2045472072204566
$ws.Range("Q8").Value = 25.86699019206522
$ws.Range("R8").Value = 232.802911728587
$ws.Range("S8").Value = 0.09727399647979741
$ws.Range("T8").Value = 0.0972739964797974

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 8.497489666666667
$ws.Range("H9").Value = 25.492469
$ws.Range("I9").Value = 0.4755576856884563
$ws.Range("J9").Value = 0.4755576856884562
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 2.355570666666666
$ws.Range("N9").Value = 7.066712
$ws.Range("O9").Value = 0.1582830603053919
$ws.Range("P9").Value = 0.1582830603053919
$ws.Range("Q9").Value = 20.01643739910311
$ws.Range("R9").Value = 180.147936591928
$ws.Range("S9").Value = 0.07527272584251854
$ws.Range("T9").Value = 0.07527272584251853

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 4.394165333333333
$ws.Range("H10").Value = 13.182496
$ws.Range("I10").Value = 0.245917227137055
$ws.Range("J10").Value = 0.2459172271370549
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 4.543446666666667
$ws.Range("N10").Value = 13.63034
$ws.Range("O10").Value = 0.305297842646339
$ws.Range("P10").Value = 0.305297842646339
$ws.Range("Q10").Value = 19.96465583651556
$ws.Range("R10").Value = 179.68190252864
$ws.Range("S10").Value = 0.07507799891451263
$ws.Range("T10").Value = 0.0750779989145126

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 4.394165333333333
$ws.Range("H11").Value = 13.182496
$ws.Range("I11").Value = 0.245917227137055
$ws.Range("J11").Value = 0.2459172271370549
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 4.938922000000001
$ws.Range("N11").Value = 14.816766
$ws.Range("O11").Value = 0.3318718898278125
$ws.Range("P11").Value = 0.3318718898278125
$ws.Range("Q11").Value = 21.70243983643734
$ws.Range("R11").Value = 195.321958527936
$ws.Range("S11").Value = 0.08161301491118986
$ws.Range("T11").Value = 0.08161301491118984

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 4.394165333333333
$ws.Range("H12").Value = 13.182496
$ws.Range("I12").Value = 0.245917227137055
$ws.Range("J12").Value = 0.2459172271370549
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 3.044074333333333
$ws.Range("N12").Value = 9.132223
$ws.Range("O12").Value = 0.2045472072204566
$ws.Range("P12").Value = 0.2045472072204566
$ws.Range("Q12").Value = 13.37616590762311
$ws.Range("R12").Value = 120.385493168608
$ws.Range("S12").Value = 0.05030168201828326
$ws.Range("T12").Value = 0.05030168201828326

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 4.394165333333333
$ws.Range("H13").Value = 13.182496
$ws.Range("I13").Value = 0.245917227137055
$ws.Range("J13").Value = 0.2459172271370549
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 2.355570666666666
$ws.Range("N13").Value = 7.066712
$ws.Range("O13").Value = 0.1582830603053919
$ws.Range("P13").Value = 0.1582830603053919
$ws.Range("Q13").Value = 10.35076696368355
$ws.Range("R13").Value = 93.156902673152
$ws.Range("S13").Value = 0.03892453129306923
$ws.Range("T13").Value = 0.03892453129306923

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 2.921228333333333
$ws.Range("H14").Value = 8.763684999999999
$ws.Range("I14").Value = 0.1634850573595927
$ws.Range("J14").Value = 0.1634850573595927
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 4.543446666666667
$ws.Range("N14").Value = 13.63034
$ws.Range("O14").Value = 0.305297842646339
$ws.Range("P14").Value = 0.305297842646339
$ws.Range("Q14").Value = 13.27244513365556
$ws.Range("R14").Value = 119.4520062029
$ws.Range("S14").Value = 0.04991163531679664
$ws.Range("T14").Value = 0.04991163531679663

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 2.921228333333333
$ws.Range("H15").Value = 8.763684999999999
$ws.Range("I15").Value = 0.1634850573595927
$ws.Range("J15").Value = 0.1634850573595927
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 4.938922000000001
$ws.Range("N15").Value = 14.816766
$ws.Range("O15").Value = 0.3318718898278125
$ws.Range("P15").Value = 0.3318718898278125
$ws.Range("Q15").Value = 14.42771888252333
$ws.Range("R15").Value = 129.84946994271
$ws.Range("S15").Value = 0.05425609494453636
$ws.Range("T15").Value = 0.05425609494453635

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 2.921228333333333
$ws.Range("H16").Value = 8.763684999999999
$ws.Range("I16").Value = 0.1634850573595927
$ws.Range("J16").Value = 0.1634850573595927
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 3.044074333333333
$ws.Range("N16").Value = 9.132223
$ws.Range("O16").Value = 0.2045472072204566
$ws.Range("P16").Value = 0.2045472072204566
$ws.Range("Q16").Value = 8.892436191306111
$ws.Range("R16").Value = 80.03192572175499
$ws.Range("S16").Value = 0.03344041190518084
$ws.Range("T16").Value = 0.03344041190518083

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 2.921228333333333
$ws.Range("H17").Value = 8.763684999999999
$ws.Range("I17").Value = 0.1634850573595927
$ws.Range("J17").Value = 0.1634850573595927
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 2.355570666666666
$ws.Range("N17").Value = 7.066712
$ws.Range("O17").Value = 0.1582830603053919
$ws.Range("P17").Value = 0.1582830603053919
$ws.Range("Q17").Value = 6.881159772635554
$ws.Range("R17").Value = 61.93043795371999
$ws.Range("S17").Value = 0.02587691519307887
$ws.Range("T17").Value = 0.02587691519307886
